# Fixed serializer, glob order lookup
# Append the newly-located Bakery order lines (rows 3-5) to the sheet.
# Values are written as text (matching the existing SKU/Qty/Cost columns,
# which are stored as strings rather than numbers) by round-tripping each
# cell through a quoted formula and then flattening it back to a plain
# value with Copy/PasteSpecial - this avoids Excel's automatic "looks like
# a number" coercion without leaving behind a new number-format style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("126160", "Peanut Butter - Bulk", "2", "65.77", "131.54"),
    @("124440", "8 Grain", "8", "68.52", "548.16"),
    @("123060", "Glaze - Donut", "1", "49.39", "49.39")
)

$startRow = 3
$r = $startRow
foreach ($rowValues in $data) {
    $c = 1
    foreach ($val in $rowValues) {
        $cell = $ws.Cells.Item($r, $c)
        $escaped = $val.Replace('"', '""')
        $cell.Formula = '="' + $escaped + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
        $c = $c + 1
    }
    $r = $r + 1
}

$excel.CutCopyMode = 0
